$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.430.72"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.603.41"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'514.05"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "'153.36"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("D9").Value = "2.612.41"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "'6.67"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "'0.345"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "'0.129"
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").Value = "3.060.47"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "60.508.30"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "'21.62"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "2.608.90"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'4.74"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "'357.47"
$ws.Range("E20").Value = "  +5.41%  "
$ws.Range("D21").Value = "'10.62"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = "  +2.42%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'61.05"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").Value = "'0.425"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("D26").Value = "2.723.85"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "0.0₃0837"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").Value = "'7.27"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "'19.42"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("D35").Value = "'150.49"
$ws.Range("E35").Value = "  -3.20%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").Value = "'0.917"
$ws.Range("E37").Value = "  +7.88%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "'0.845"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").Value = "'36.22"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").Value = "'3.74"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'288.01"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'0.101"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").Value = "'0.619"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'0.997"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "'0.0555"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("D48").Value = "'19.64"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "'4.95"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").Value = "'0.0236"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("E51").Value = "  +0.42%  "
